$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: change "A" (shift marker) to numeric shift numbers
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 2

# Move the "**" note (and its formatting) from C6 to A6
$null = $ws.Range("C6").Copy()
$null = $ws.Range("A6").PasteSpecial(-4122)
$ws.Range("A6").Value = "**"
$null = $ws.Range("C6").Clear()

# Update the active selection to match the saved view state
$null = $ws.Range("D11").Select()
